$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.7246536666666668
$ws.Range("H2").Value = 2.173961
$ws.Range("I2").Value = 0.1791272621505297
$ws.Range("J2").Value = 0.1791272621505298
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.024839333333334
$ws.Range("N2").Value = 9.074518000000001
$ws.Range("O2").Value = 0.1801507982970389
$ws.Range("P2").Value = 0.1801507982970388
$ws.Range("Q2").Value = 2.191960913977556
$ws.Range("R2").Value = 19.72764822579801
$ws.Range("S2").Value = 0.03226991927318089
$ws.Range("T2").Value = 0.03226991927318089

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.7246536666666668
$ws.Range("H3").Value = 2.173961
$ws.Range("I3").Value = 0.1791272621505297
$ws.Range("J3").Value = 0.1791272621505298
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.626140333333333
$ws.Range("N3").Value = 4.878420999999999
$ws.Range("O3").Value = 0.09684827751501936
$ws.Range("P3").Value = 0.09684827751501934
$ws.Range("Q3").Value = 1.178388555064555
$ws.Range("R3").Value = 10.605496995581
$ws.Range("S3").Value = 0.01734816679526013
$ws.Range("T3").Value = 0.01734816679526013

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.7246536666666668
$ws.Range("H4").Value = 2.173961
$ws.Range("I4").Value = 0.1791272621505297
$ws.Range("J4").Value = 0.1791272621505298
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.640628666666666
$ws.Range("N4").Value = 13.921886
$ws.Range("O4").Value = 0.2763825997921178
$ws.Range("P4").Value = 0.2763825997921177
$ws.Range("Q4").Value = 3.362848578938444
$ws.Range("R4").Value = 30.265637210446
$ws.Range("S4").Value = 0.04950765840680763
$ws.Range("T4").Value = 0.04950765840680763

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.7246536666666668
$ws.Range("H5").Value = 2.173961
$ws.Range("I5").Value = 0.1791272621505297
$ws.Range("J5").Value = 0.1791272621505298
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.498988000000001
$ws.Range("N5").Value = 22.496964
$ws.Range("O5").Value = 0.4466183243958241
$ws.Range("P5").Value = 0.446618324395824
$ws.Range("Q5").Value = 5.434169150489335
$ws.Range("R5").Value = 48.90752235440401
$ws.Range("S5").Value = 0.0800015176752811
$ws.Range("T5").Value = 0.0800015176752811

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.320814666666667
$ws.Range("H6").Value = 9.962444
$ws.Range("I6").Value = 0.8208727378494701
$ws.Range("J6").Value = 0.8208727378494702
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.024839333333334
$ws.Range("N6").Value = 9.074518000000001
$ws.Range("O6").Value = 0.1801507982970389
$ws.Range("P6").Value = 0.1801507982970388
$ws.Range("Q6").Value = 10.04493082244356
$ws.Range("R6").Value = 90.40437740199201
$ws.Range("S6").Value = 0.147880879023858
$ws.Range("T6").Value = 0.147880879023858

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.320814666666667
$ws.Range("H7").Value = 9.962444
$ws.Range("I7").Value = 0.8208727378494701
$ws.Range("J7").Value = 0.8208727378494702
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.626140333333333
$ws.Range("N7").Value = 4.878420999999999
$ws.Range("O7").Value = 0.09684827751501936
$ws.Range("P7").Value = 0.09684827751501934
$ws.Range("Q7").Value = 5.400110668991554
$ws.Range("R7").Value = 48.60099602092399
$ws.Range("S7").Value = 0.07950011071975922
$ws.Range("T7").Value = 0.07950011071975921

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.320814666666667
$ws.Range("H8").Value = 9.962444
$ws.Range("I8").Value = 0.8208727378494701
$ws.Range("J8").Value = 0.8208727378494702
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.640628666666666
$ws.Range("N8").Value = 13.921886
$ws.Range("O8").Value = 0.2763825997921178
$ws.Range("P8").Value = 0.2763825997921177
$ws.Range("Q8").Value = 15.41066773882044
$ws.Range("R8").Value = 138.696009649384
$ws.Range("S8").Value = 0.2268749413853101
$ws.Range("T8").Value = 0.2268749413853101

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.320814666666667
$ws.Range("H9").Value = 9.962444
$ws.Range("I9").Value = 0.8208727378494701
$ws.Range("J9").Value = 0.8208727378494702
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.498988000000001
$ws.Range("N9").Value = 22.496964
$ws.Range("O9").Value = 0.4466183243958241
$ws.Range("P9").Value = 0.446618324395824
$ws.Range("Q9").Value = 24.90274933555733
$ws.Range("R9").Value = 224.124744020016
$ws.Range("S9").Value = 0.3666168067205429
$ws.Range("T9").Value = 0.3666168067205429
